# "Generate Report for Handback"
#
# The localization-status report is refreshed after a handback completes:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this text lives in every Status-ish cell on Overview + each language sheet)
#   - Each language sheet's per-file rows gain their handback details:
#       * "Latest Target File"    (J) -> link to the source .md (same as column A)
#       * "Latest Handback File"  (K) -> the handback xlf file name (same base name
#                                         as the "Latest Handoff File" in column G)
#       * "Latest Handback DateTime" (L) -> the timestamp the handback finished
#
$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70442f9be640adc2b130e59bc9e81f2090c14720/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet: both language status columns (zh-cn / de-de) move to the
# "handed back" state for both tracked files.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------------
# Per-language detail sheets: zh-cn finished handback at 16:54:09,
# de-de finished handback at 16:54:28.
# ---------------------------------------------------------------------------
$languageSheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-12-07 16:54:09" },
    @{ Name = "de-de"; HandbackTime = "2016-12-07 16:54:28" }
)

foreach ($lang in $languageSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Widen columns to fit the new, longer handback file names/links.
    $ws.Columns.Item(3).ColumnWidth = 29.9777050018311
    $ws.Columns.Item(10).ColumnWidth = 40
    $ws.Columns.Item(11).ColumnWidth = 40

    for ($row = 2; $row -le 3; $row++) {
        # Status column moves to "handed back".
        $ws.Cells.Item($row, 3).Value = $newStatus

        # Latest Target File (J): hyperlink to the same source .md as column A.
        $sourceName = $ws.Cells.Item($row, 1).Value()
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 10), ($repoBase + $sourceName), "", "", $sourceName)

        # Latest Handback File (K): the handback xlf, same base name as the
        # handoff xlf already recorded in column G.
        $handoffFile = $ws.Cells.Item($row, 7).Value()
        $ws.Cells.Item($row, 11).Value = $handoffFile

        # Latest Handback DateTime (L)
        $ws.Cells.Item($row, 12).Value = $lang.HandbackTime
    }
}
